# "Generate Report for Handback" - refresh the handoff/handback timestamps
# that the report generator stamps into the handback-status workbook.
$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview!G2 "Latest HO Xliff Generate Date" == the newest per-locale
# "Correspond Handoff Datetime" (de-de was generated last), so both cells
# carry the same refreshed timestamp.
$wsOverview.Range("G2").Value = "2016-11-15 16:18:19"
$wsDeDe.Range("H2").Value = "2016-11-15 16:18:19"

# zh-cn row: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H2").Value = "2016-11-15 16:18:04"
$wsZhCn.Range("K2").Value = "2016-11-15 16:19:01"

# de-de row: Correspond Handback DateTime
$wsDeDe.Range("K2").Value = "2016-11-15 16:19:19"
